$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 19.163986
$ws.Cells.Item(2, 8).Value = 57.491958
$ws.Cells.Item(2, 9).Value = 0.1197574615923936
$ws.Cells.Item(2, 10).Value = 0.1197574615923936
$ws.Cells.Item(2, 13).Value = 68.63737500000001
$ws.Cells.Item(2, 14).Value = 205.912125
$ws.Cells.Item(2, 15).Value = 0.5415701538216162
$ws.Cells.Item(2, 16).Value = 0.5415701538216162
$ws.Cells.Item(2, 17).Value = 1315.36569357675
$ws.Cells.Item(2, 18).Value = 11838.29124219075
$ws.Cells.Item(2, 19).Value = 0.06485706689587889
$ws.Cells.Item(2, 20).Value = 0.0648570668958789
$ws.Cells.Item(3, 7).Value = 19.163986
$ws.Cells.Item(3, 8).Value = 57.491958
$ws.Cells.Item(3, 9).Value = 0.1197574615923936
$ws.Cells.Item(3, 10).Value = 0.1197574615923936
$ws.Cells.Item(3, 15).Value = 0.08718851262838957
$ws.Cells.Item(3, 16).Value = 0.08718851262838957
$ws.Cells.Item(3, 17).Value = 211.7634762109547
$ws.Cells.Item(3, 18).Value = 1905.871285898592
$ws.Cells.Item(3, 19).Value = 0.01044147495239229
$ws.Cells.Item(3, 20).Value = 0.01044147495239229
$ws.Cells.Item(4, 7).Value = 19.163986
$ws.Cells.Item(4, 8).Value = 57.491958
$ws.Cells.Item(4, 9).Value = 0.1197574615923936
$ws.Cells.Item(4, 10).Value = 0.1197574615923936
$ws.Cells.Item(4, 13).Value = 16.21089566666667
$ws.Cells.Item(4, 14).Value = 48.632687
$ws.Cells.Item(4, 15).Value = 0.1279089892319285
$ws.Cells.Item(4, 16).Value = 0.1279089892319285
$ws.Cells.Item(4, 17).Value = 310.6653776034607
$ws.Cells.Item(4, 18).Value = 2795.988398431146
$ws.Cells.Item(4, 19).Value = 0.01531805586526456
$ws.Cells.Item(4, 20).Value = 0.01531805586526457
$ws.Cells.Item(5, 7).Value = 19.163986
$ws.Cells.Item(5, 8).Value = 57.491958
$ws.Cells.Item(5, 9).Value = 0.1197574615923936
$ws.Cells.Item(5, 10).Value = 0.1197574615923936
$ws.Cells.Item(5, 13).Value = 20.32546233333333
$ws.Cells.Item(5, 14).Value = 60.976387
$ws.Cells.Item(5, 15).Value = 0.1603741949973873
$ws.Cells.Item(5, 16).Value = 0.1603741949973873
$ws.Cells.Item(5, 17).Value = 389.5168755995274
$ws.Cells.Item(5, 18).Value = 3505.651880395747
$ws.Cells.Item(5, 19).Value = 0.01920600649781065
$ws.Cells.Item(5, 20).Value = 0.01920600649781065
$ws.Cells.Item(6, 7).Value = 19.163986
$ws.Cells.Item(6, 8).Value = 57.491958
$ws.Cells.Item(6, 9).Value = 0.1197574615923936
$ws.Cells.Item(6, 10).Value = 0.1197574615923936
$ws.Cells.Item(6, 13).Value = 10.513928
$ws.Cells.Item(6, 14).Value = 31.541784
$ws.Cells.Item(6, 15).Value = 0.08295814932067838
$ws.Cells.Item(6, 16).Value = 0.08295814932067838
$ws.Cells.Item(6, 17).Value = 201.488768997008
$ws.Cells.Item(6, 18).Value = 1813.398920973072
$ws.Cells.Item(6, 19).Value = 0.009934857381047194
$ws.Cells.Item(6, 20).Value = 0.009934857381047194
$ws.Cells.Item(7, 9).Value = 0.150345281456851
$ws.Cells.Item(7, 10).Value = 0.1503452814568511
$ws.Cells.Item(7, 13).Value = 68.63737500000001
$ws.Cells.Item(7, 14).Value = 205.912125
$ws.Cells.Item(7, 15).Value = 0.5415701538216162
$ws.Cells.Item(7, 16).Value = 0.5415701538216162
$ws.Cells.Item(7, 17).Value = 1651.329468660375
$ws.Cells.Item(7, 18).Value = 14861.96521794337
$ws.Cells.Item(7, 19).Value = 0.08142251720494099
$ws.Cells.Item(7, 20).Value = 0.08142251720494101
$ws.Cells.Item(8, 9).Value = 0.150345281456851
$ws.Cells.Item(8, 10).Value = 0.1503452814568511
$ws.Cells.Item(8, 15).Value = 0.08718851262838957
$ws.Cells.Item(8, 16).Value = 0.08718851262838957
$ws.Cells.Item(8, 19).Value = 0.01310838147091944
$ws.Cells.Item(8, 20).Value = 0.01310838147091944
$ws.Cells.Item(9, 9).Value = 0.150345281456851
$ws.Cells.Item(9, 10).Value = 0.1503452814568511
$ws.Cells.Item(9, 13).Value = 16.21089566666667
$ws.Cells.Item(9, 14).Value = 48.632687
$ws.Cells.Item(9, 15).Value = 0.1279089892319285
$ws.Cells.Item(9, 16).Value = 0.1279089892319285
$ws.Cells.Item(9, 17).Value = 390.0138915240486
$ws.Cells.Item(9, 18).Value = 3510.125023716437
$ws.Cells.Item(9, 19).Value = 0.01923051298693562
$ws.Cells.Item(9, 20).Value = 0.01923051298693562
$ws.Cells.Item(10, 9).Value = 0.150345281456851
$ws.Cells.Item(10, 10).Value = 0.1503452814568511
$ws.Cells.Item(10, 13).Value = 20.32546233333333
$ws.Cells.Item(10, 14).Value = 60.976387
$ws.Cells.Item(10, 15).Value = 0.1603741949973873
$ws.Cells.Item(10, 16).Value = 0.1603741949973873
$ws.Cells.Item(10, 17).Value = 489.0052236872374
$ws.Cells.Item(10, 18).Value = 4401.047013185137
$ws.Cells.Item(10, 19).Value = 0.02411150348529811
$ws.Cells.Item(10, 20).Value = 0.02411150348529811
$ws.Cells.Item(11, 9).Value = 0.150345281456851
$ws.Cells.Item(11, 10).Value = 0.1503452814568511
$ws.Cells.Item(11, 13).Value = 10.513928
$ws.Cells.Item(11, 14).Value = 31.541784
$ws.Cells.Item(11, 15).Value = 0.08295814932067838
$ws.Cells.Item(11, 16).Value = 0.08295814932067838
$ws.Cells.Item(11, 17).Value = 252.9519687746427
$ws.Cells.Item(11, 18).Value = 2276.567718971784
$ws.Cells.Item(11, 19).Value = 0.01247236630875687
$ws.Cells.Item(11, 20).Value = 0.01247236630875687
$ws.Cells.Item(12, 7).Value = 61.341815
$ws.Cells.Item(12, 8).Value = 184.025445
$ws.Cells.Item(12, 9).Value = 0.3833304853108436
$ws.Cells.Item(12, 10).Value = 0.3833304853108436
$ws.Cells.Item(12, 13).Value = 68.63737500000001
$ws.Cells.Item(12, 14).Value = 205.912125
$ws.Cells.Item(12, 15).Value = 0.5415701538216162
$ws.Cells.Item(12, 16).Value = 0.5415701538216162
$ws.Cells.Item(12, 17).Value = 4210.341159335625
$ws.Cells.Item(12, 18).Value = 37893.07043402062
$ws.Cells.Item(12, 19).Value = 0.2076003498943083
$ws.Cells.Item(12, 20).Value = 0.2076003498943084
$ws.Cells.Item(13, 7).Value = 61.341815
$ws.Cells.Item(13, 8).Value = 184.025445
$ws.Cells.Item(13, 9).Value = 0.3833304853108436
$ws.Cells.Item(13, 10).Value = 0.3833304853108436
$ws.Cells.Item(13, 15).Value = 0.08718851262838957
$ws.Cells.Item(13, 16).Value = 0.08718851262838957
$ws.Cells.Item(13, 17).Value = 677.8316359388533
$ws.Cells.Item(13, 18).Value = 6100.48472344968
$ws.Cells.Item(13, 19).Value = 0.03342201485937119
$ws.Cells.Item(13, 20).Value = 0.03342201485937119
$ws.Cells.Item(14, 7).Value = 61.341815
$ws.Cells.Item(14, 8).Value = 184.025445
$ws.Cells.Item(14, 9).Value = 0.3833304853108436
$ws.Cells.Item(14, 10).Value = 0.3833304853108436
$ws.Cells.Item(14, 13).Value = 16.21089566666667
$ws.Cells.Item(14, 14).Value = 48.632687
$ws.Cells.Item(14, 15).Value = 0.1279089892319285
$ws.Cells.Item(14, 16).Value = 0.1279089892319285
$ws.Cells.Item(14, 17).Value = 994.4057629689685
$ws.Cells.Item(14, 18).Value = 8949.651866720715
$ws.Cells.Item(14, 19).Value = 0.04903141491789462
$ws.Cells.Item(14, 20).Value = 0.04903141491789462
$ws.Cells.Item(15, 7).Value = 61.341815
$ws.Cells.Item(15, 8).Value = 184.025445
$ws.Cells.Item(15, 9).Value = 0.3833304853108436
$ws.Cells.Item(15, 10).Value = 0.3833304853108436
$ws.Cells.Item(15, 13).Value = 20.32546233333333
$ws.Cells.Item(15, 14).Value = 60.976387
$ws.Cells.Item(15, 15).Value = 0.1603741949973873
$ws.Cells.Item(15, 16).Value = 0.1603741949973873
$ws.Cells.Item(15, 17).Value = 1246.800750240802
$ws.Cells.Item(15, 18).Value = 11221.20675216721
$ws.Cells.Item(15, 19).Value = 0.06147631799968435
$ws.Cells.Item(15, 20).Value = 0.06147631799968436
$ws.Cells.Item(16, 7).Value = 61.341815
$ws.Cells.Item(16, 8).Value = 184.025445
$ws.Cells.Item(16, 9).Value = 0.3833304853108436
$ws.Cells.Item(16, 10).Value = 0.3833304853108436
$ws.Cells.Item(16, 13).Value = 10.513928
$ws.Cells.Item(16, 14).Value = 31.541784
$ws.Cells.Item(16, 15).Value = 0.08295814932067838
$ws.Cells.Item(16, 16).Value = 0.08295814932067838
$ws.Cells.Item(16, 17).Value = 644.94342629932
$ws.Cells.Item(16, 18).Value = 5804.49083669388
$ws.Cells.Item(16, 19).Value = 0.03180038763958507
$ws.Cells.Item(16, 20).Value = 0.03180038763958508
$ws.Cells.Item(17, 7).Value = 7.095824666666666
$ws.Cells.Item(17, 8).Value = 21.287474
$ws.Cells.Item(17, 9).Value = 0.04434244264135302
$ws.Cells.Item(17, 10).Value = 0.04434244264135302
$ws.Cells.Item(17, 13).Value = 68.63737500000001
$ws.Cells.Item(17, 14).Value = 205.912125
$ws.Cells.Item(17, 15).Value = 0.5415701538216162
$ws.Cells.Item(17, 16).Value = 0.5415701538216162
$ws.Cells.Item(17, 17).Value = 487.03877858025
$ws.Cells.Item(17, 18).Value = 4383.34900722225
$ws.Cells.Item(17, 19).Value = 0.02401454348210375
$ws.Cells.Item(17, 20).Value = 0.02401454348210375
$ws.Cells.Item(18, 7).Value = 7.095824666666666
$ws.Cells.Item(18, 8).Value = 21.287474
$ws.Cells.Item(18, 9).Value = 0.04434244264135302
$ws.Cells.Item(18, 10).Value = 0.04434244264135302
$ws.Cells.Item(18, 15).Value = 0.08718851262838957
$ws.Cells.Item(18, 16).Value = 0.08718851262838957
$ws.Cells.Item(18, 17).Value = 78.40939238824177
$ws.Cells.Item(18, 18).Value = 705.684531494176
$ws.Cells.Item(18, 19).Value = 0.003866151620209248
$ws.Cells.Item(18, 20).Value = 0.003866151620209248
$ws.Cells.Item(19, 7).Value = 7.095824666666666
$ws.Cells.Item(19, 8).Value = 21.287474
$ws.Cells.Item(19, 9).Value = 0.04434244264135302
$ws.Cells.Item(19, 10).Value = 0.04434244264135302
$ws.Cells.Item(19, 13).Value = 16.21089566666667
$ws.Cells.Item(19, 14).Value = 48.632687
$ws.Cells.Item(19, 15).Value = 0.1279089892319285
$ws.Cells.Item(19, 16).Value = 0.1279089892319285
$ws.Cells.Item(19, 17).Value = 115.0296733402931
$ws.Cells.Item(19, 18).Value = 1035.267060062638
$ws.Cells.Item(19, 19).Value = 0.005671797018330231
$ws.Cells.Item(19, 20).Value = 0.005671797018330231
$ws.Cells.Item(20, 7).Value = 7.095824666666666
$ws.Cells.Item(20, 8).Value = 21.287474
$ws.Cells.Item(20, 9).Value = 0.04434244264135302
$ws.Cells.Item(20, 10).Value = 0.04434244264135302
$ws.Cells.Item(20, 13).Value = 20.32546233333333
$ws.Cells.Item(20, 14).Value = 60.976387
$ws.Cells.Item(20, 15).Value = 0.1603741949973873
$ws.Cells.Item(20, 16).Value = 0.1603741949973873
$ws.Cells.Item(20, 17).Value = 144.2259169862709
$ws.Cells.Item(20, 18).Value = 1298.033252876438
$ws.Cells.Item(20, 19).Value = 0.007111383542824813
$ws.Cells.Item(20, 20).Value = 0.007111383542824813
$ws.Cells.Item(21, 7).Value = 7.095824666666666
$ws.Cells.Item(21, 8).Value = 21.287474
$ws.Cells.Item(21, 9).Value = 0.04434244264135302
$ws.Cells.Item(21, 10).Value = 0.04434244264135302
$ws.Cells.Item(21, 13).Value = 10.513928
$ws.Cells.Item(21, 14).Value = 31.541784
$ws.Cells.Item(21, 15).Value = 0.08295814932067838
$ws.Cells.Item(21, 16).Value = 0.08295814932067838
$ws.Cells.Item(21, 17).Value = 74.60498964595733
$ws.Cells.Item(21, 18).Value = 671.4449068136159
$ws.Cells.Item(21, 19).Value = 0.00367856697788498
$ws.Cells.Item(21, 20).Value = 0.00367856697788498
$ws.Cells.Item(22, 7).Value = 48.362939
$ws.Cells.Item(22, 8).Value = 145.088817
$ws.Cells.Item(22, 9).Value = 0.3022243289985588
$ws.Cells.Item(22, 10).Value = 0.3022243289985588
$ws.Cells.Item(22, 13).Value = 68.63737500000001
$ws.Cells.Item(22, 14).Value = 205.912125
$ws.Cells.Item(22, 15).Value = 0.5415701538216162
$ws.Cells.Item(22, 16).Value = 0.5415701538216162
$ws.Cells.Item(22, 17).Value = 3319.505180245126
$ws.Cells.Item(22, 18).Value = 29875.54662220613
$ws.Cells.Item(22, 19).Value = 0.1636756763443842
$ws.Cells.Item(22, 20).Value = 0.1636756763443842
$ws.Cells.Item(23, 7).Value = 48.362939
$ws.Cells.Item(23, 8).Value = 145.088817
$ws.Cells.Item(23, 9).Value = 0.3022243289985588
$ws.Cells.Item(23, 10).Value = 0.3022243289985588
$ws.Cells.Item(23, 15).Value = 0.08718851262838957
$ws.Cells.Item(23, 16).Value = 0.08718851262838957
$ws.Cells.Item(23, 17).Value = 534.4140870494454
$ws.Cells.Item(23, 18).Value = 4809.726783445009
$ws.Cells.Item(23, 19).Value = 0.02635048972549741
$ws.Cells.Item(23, 20).Value = 0.02635048972549741
$ws.Cells.Item(24, 7).Value = 48.362939
$ws.Cells.Item(24, 8).Value = 145.088817
$ws.Cells.Item(24, 9).Value = 0.3022243289985588
$ws.Cells.Item(24, 10).Value = 0.3022243289985588
$ws.Cells.Item(24, 13).Value = 16.21089566666667
$ws.Cells.Item(24, 14).Value = 48.632687
$ws.Cells.Item(24, 15).Value = 0.1279089892319285
$ws.Cells.Item(24, 16).Value = 0.1279089892319285
$ws.Cells.Item(24, 17).Value = 784.0065582623645
$ws.Cells.Item(24, 18).Value = 7056.05902436128
$ws.Cells.Item(24, 19).Value = 0.03865720844350348
$ws.Cells.Item(24, 20).Value = 0.03865720844350348
$ws.Cells.Item(25, 7).Value = 48.362939
$ws.Cells.Item(25, 8).Value = 145.088817
$ws.Cells.Item(25, 9).Value = 0.3022243289985588
$ws.Cells.Item(25, 10).Value = 0.3022243289985588
$ws.Cells.Item(25, 13).Value = 20.32546233333333
$ws.Cells.Item(25, 14).Value = 60.976387
$ws.Cells.Item(25, 15).Value = 0.1603741949973873
$ws.Cells.Item(25, 16).Value = 0.1603741949973873
$ws.Cells.Item(25, 17).Value = 982.9990949737978
$ws.Cells.Item(25, 18).Value = 8846.991854764179
$ws.Cells.Item(25, 19).Value = 0.04846898347176941
$ws.Cells.Item(25, 20).Value = 0.04846898347176941
$ws.Cells.Item(26, 7).Value = 48.362939
$ws.Cells.Item(26, 8).Value = 145.088817
$ws.Cells.Item(26, 9).Value = 0.3022243289985588
$ws.Cells.Item(26, 10).Value = 0.3022243289985588
$ws.Cells.Item(26, 13).Value = 10.513928
$ws.Cells.Item(26, 14).Value = 31.541784
$ws.Cells.Item(26, 15).Value = 0.08295814932067838
$ws.Cells.Item(26, 16).Value = 0.08295814932067838
$ws.Cells.Item(26, 17).Value = 508.484458514392
$ws.Cells.Item(26, 18).Value = 4576.360126629528
$ws.Cells.Item(26, 19).Value = 0.02507197101340427
$ws.Cells.Item(26, 20).Value = 0.02507197101340427
